$d = $word.ActiveDocument

# The "Profile" paragraph originally reads:
#   "A Computer Science student on track to graduate in four months, I'm deeply
#    passionate about using advanced technology to solve complex problems. My
#    academic journey has equipped me with a solid grasp of computer science
#    principles, various programming languages, and software development
#    methodologies. I've also gained invaluable hands-on experience through
#    internship at prominent tech companies like Tech Mahindra."
#
# It becomes:
#   "A Computer Science student on track to graduate in 4 months, I am deeply
#    passionate about using advanced technology to solve complex problems. My
#    academic journey has equipped me with a solid grasp of computer science
#    principles, various programming languages, and software development
#    methodologies. I have also gained invaluable hands-on experience through
#    internship at prominent tech companies like Tech Mahindra."
#
# i.e. "four" -> "4", "I'm" -> "I am", "I've" -> "I have".

$found1 = $d.Content.Find.Execute("four months", $true, $false, $false, $false, $false, $true, 1, $false, "4 months", 2)
Write-Output "replaced 'four months': $found1"

$found2 = $d.Content.Find.Execute("I'm deeply", $true, $false, $false, $false, $false, $true, 1, $false, "I am deeply", 2)
Write-Output "replaced ""I'm deeply"": $found2"

$found3 = $d.Content.Find.Execute("I've also", $true, $false, $false, $false, $false, $true, 1, $false, "I have also", 2)
Write-Output "replaced ""I've also"": $found3"
